$wb = $excel.ActiveWorkbook

# --- Logs sheet: append row 11 with the new mail log entry ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A11").Value = "Klacht over levering"
$logs.Range("B11").Value = "klantenservice@testbedrijf123.nl"
$logs.Range("D11").Value = "Intern verzoek / Actie voor medewerker"
$logs.Range("F11").Value = "2025-08-19 19:50:03"
$logs.Range("G11").Value = "Nee"
$logs.Range("H11").Value = "Ja"
$logs.Range("I11").Value = "Nee"
$logs.Range("J11").Value = "Nee"

# Extend the conditional formatting ranges to include the new row 11
$dFcs = $logs.Range("D2:D11").FormatConditions
for ($i = 1; $i -le $dFcs.Count; $i++) {
    $dFcs.Item($i).ModifyAppliesToRange($logs.Range("D2:D11"))
}

$gFcs = $logs.Range("G2:G11").FormatConditions
for ($i = 1; $i -le $gFcs.Count; $i++) {
    $gFcs.Item($i).ModifyAppliesToRange($logs.Range("G2:G11"))
}

$hFcs = $logs.Range("H2:H11").FormatConditions
for ($i = 1; $i -le $hFcs.Count; $i++) {
    $hFcs.Item($i).ModifyAppliesToRange($logs.Range("H2:H11"))
}

$iFcs = $logs.Range("I2:I11").FormatConditions
for ($i = 1; $i -le $iFcs.Count; $i++) {
    $iFcs.Item($i).ModifyAppliesToRange($logs.Range("I2:I11"))
}

$jFcs = $logs.Range("J2:J11").FormatConditions
for ($i = 1; $i -le $jFcs.Count; $i++) {
    $jFcs.Item($i).ModifyAppliesToRange($logs.Range("J2:J11"))
}

# --- Dashboard sheet: update the count for the category ---
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 10
